# Add a "Save" column (H) to the sheet, mirroring the header style used by
# the existing header row (column G uses the bold/bordered/centered style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: copy formatting from G1 (reuses the existing header style)
# then set its value.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Data cells for the new "Save" column.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
